$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '291.53'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-3.23%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '30.65'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-6.26%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '4.951'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '0.29%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07215'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-6.69%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.815'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '-8.11%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '7.690'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '-1.83%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.761'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-0.98%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8979'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-2.44%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1663'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-5.41%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07708'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-0.96%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08045'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-6.39%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03037'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-4.25%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.1002'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-0.06%'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '-1.20%'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.005727'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-2.75%'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.470'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '0.27%'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.083'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-3.28%'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-0.87%'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1289'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '-2.87%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.049'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-5.25%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.2252'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '13.03%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.04505'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-0.94%'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '-0.77%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004012'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-9.05%'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '-0.09%'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '-5.84%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04419'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '-5.59%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007309'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-6.26%'
$ws.Range('B42').Value = 'Dexo'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.009922'
$ws.Range('E42').Value = '--%'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1308'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-3.19%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.002008'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-13.56%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.009515'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-16.81%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00005964'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-4.46%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-0.03%'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.246'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '173.66%'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.003002'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '-3.31%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '-0.03%'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '-0.03%'
